$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.92074263832864
$ws.Range("D2").Value = 9.667037446485168
$ws.Range("E2").Value = 13.75795227573975
$ws.Range("F2").Value = 30.96085758376382
$ws.Range("G2").Value = 33.36414497991095
$ws.Range("H2").Value = 14.37340777224149
$ws.Range("J2").Value = 9.73228269459122
$ws.Range("N2").Value = 18.99769123883711
$ws.Range("O2").Value = 22.74505173743382

$ws.Range("C3").Value = 4.749564103841492
$ws.Range("D3").Value = 9.673078951589259
$ws.Range("E3").Value = 13.72986351928126
$ws.Range("F3").Value = 30.60052972713536
$ws.Range("G3").Value = 32.53886467159878
$ws.Range("H3").Value = 14.3199452415204
$ws.Range("J3").Value = 9.734587955838204
$ws.Range("N3").Value = 18.40031508502701
$ws.Range("O3").Value = 22.51172477890321

$ws.Range("C4").Value = 4.642857908240394
$ws.Range("D4").Value = 9.678311898181713
$ws.Range("E4").Value = 13.71559754511369
$ws.Range("F4").Value = 30.38616264653965
$ws.Range("G4").Value = 32.03238208946642
$ws.Range("H4").Value = 14.29036715736524
$ws.Range("J4").Value = 9.737749318922146
$ws.Range("N4").Value = 18.02485520896362
$ws.Range("O4").Value = 22.3736556215559

$ws.Range("C5").Value = 4.599049569503514
$ws.Range("D5").Value = 9.680827631217131
$ws.Range("E5").Value = 13.71053700801886
$ws.Range("F5").Value = 30.30063607324135
$ws.Range("G5").Value = 31.82639236883763
$ws.Range("H5").Value = 14.27913879270429
$ws.Range("J5").Value = 9.739476756589896
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 22.31875788180941

$ws.Range("C6").Value = 4.591758020446908
$ws.Range("D6").Value = 9.681268518400826
$ws.Range("E6").Value = 13.70974227467605
$ws.Range("F6").Value = 30.28654782818233
$ws.Range("G6").Value = 31.79222273498753
$ws.Range("H6").Value = 14.27732437483244
$ws.Range("J6").Value = 9.739790120798029
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 22.30972639537129

$ws.Range("C7").Value = 4.642268300577506
$ws.Range("D7").Value = 9.678344274245845
$ws.Range("E7").Value = 13.71552624415191
$ws.Range("F7").Value = 30.38500166663774
$ws.Range("G7").Value = 32.02960193277746
$ws.Range("H7").Value = 14.29021237756461
$ws.Range("J7").Value = 9.737770837578809
$ws.Range("N7").Value = 18.02277304767603
$ws.Range("O7").Value = 22.37290964086522

$ws.Range("C8").Value = 4.862097684253743
$ws.Range("D8").Value = 9.668804505872741
$ws.Range("E8").Value = 13.74765091334313
$ws.Range("F8").Value = 30.83525117583891
$ws.Range("G8").Value = 33.07974150062998
$ws.Range("H8").Value = 14.35430495193129
$ws.Range("J8").Value = 9.732715306555381
$ws.Range("N8").Value = 18.79364780656867
$ws.Range("O8").Value = 22.66355957609677

$ws.Range("C9").Value = 5.277361014224199
$ws.Range("D9").Value = 9.662173792963561
$ws.Range("E9").Value = 13.83411694412132
$ws.Range("F9").Value = 31.76813656619932
$ws.Range("G9").Value = 35.12644522169612
$ws.Range("H9").Value = 14.50535820220201
$ws.Range("J9").Value = 9.736643243213134
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 23.27195346769959

$ws.Range("C10").Value = 5.569097360694976
$ws.Range("D10").Value = 9.664645213267283
$ws.Range("E10").Value = 13.91166879830033
$ws.Range("F10").Value = 32.47761986149001
$ws.Range("G10").Value = 36.60406146725441
$ws.Range("H10").Value = 14.63119705848085
$ws.Range("J10").Value = 9.747943953925418
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 23.73848429676603

$ws.Range("C11").Value = 5.698257383232274
$ws.Range("D11").Value = 9.66735770262812
$ws.Range("E11").Value = 13.94992080461271
$ws.Range("F11").Value = 32.80423388496732
$ws.Range("G11").Value = 37.26706215651713
$ws.Range("H11").Value = 14.69152319190206
$ws.Range("J11").Value = 9.754903559497093
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 23.95410895876936

$ws.Range("C12").Value = 5.746609374124008
$ws.Range("D12").Value = 9.668612479726175
$ws.Range("E12").Value = 13.96482610652037
$ws.Range("F12").Value = 32.928354529016
$ws.Range("G12").Value = 37.51653606913369
$ws.Range("H12").Value = 14.71479642540435
$ws.Range("J12").Value = 9.757799456284602
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 24.03617626535241

$ws.Range("C13").Value = 5.736221423675627
$ws.Range("D13").Value = 9.668332131241389
$ws.Range("E13").Value = 13.96159742434225
$ws.Range("F13").Value = 32.90160525344557
$ws.Range("G13").Value = 37.46288220902905
$ws.Range("H13").Value = 14.70976526868489
$ws.Range("J13").Value = 9.757164209302486
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 24.01848431984365

$ws.Range("C14").Value = 5.702246767003802
$ws.Range("D14").Value = 9.667456376952487
$ws.Range("E14").Value = 13.95113869827613
$ws.Range("F14").Value = 32.81443723296704
$ws.Range("G14").Value = 37.28761987037135
$ws.Range("H14").Value = 14.69342937958439
$ws.Range("J14").Value = 9.755136593949604
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 23.96085279118066

$ws.Range("C15").Value = 5.681362318579184
$ws.Range("D15").Value = 9.666949569670352
$ws.Range("E15").Value = 13.94478690389393
$ws.Range("F15").Value = 32.7610980207413
$ws.Range("G15").Value = 37.18005186979817
$ws.Range("H15").Value = 14.68347861846534
$ws.Range("J15").Value = 9.75392850268393
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 23.92560366731023

$ws.Range("C16").Value = 5.560580797131255
$ws.Range("D16").Value = 9.664499841887622
$ws.Range("E16").Value = 13.90922814636434
$ws.Range("F16").Value = 32.45634264733502
$ws.Range("G16").Value = 36.56052458280225
$ws.Range("H16").Value = 14.6273154247838
$ws.Range("J16").Value = 9.747525637272394
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 23.72445479040506

$ws.Range("C17").Value = 5.485539969987388
$ws.Range("D17").Value = 9.663403286992464
$ws.Range("E17").Value = 13.88817041975878
$ws.Range("F17").Value = 32.27029187686021
$ws.Range("G17").Value = 36.17791568812812
$ws.Range("H17").Value = 14.59364058654374
$ws.Range("J17").Value = 9.744062813389085
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 23.60187428054179

$ws.Range("C18").Value = 5.442046567105474
$ws.Range("D18").Value = 9.662922122700738
$ws.Range("E18").Value = 13.87633891779273
$ws.Range("F18").Value = 32.16365129537094
$ws.Range("G18").Value = 35.95700053572269
$ws.Range("H18").Value = 14.57456264014793
$ws.Range("J18").Value = 9.742242366009144
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 23.53169342163403

$ws.Range("C19").Value = 5.427264976841398
$ws.Range("D19").Value = 9.662784914726263
$ws.Range("E19").Value = 13.87238134116863
$ws.Range("F19").Value = 32.12761190950909
$ws.Range("G19").Value = 35.88206529527492
$ws.Range("H19").Value = 14.56815355859889
$ws.Range("J19").Value = 9.741655441399654
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 23.50798931408281

$ws.Range("C20").Value = 5.493562905250962
$ws.Range("D20").Value = 9.663504544996403
$ws.Range("E20").Value = 13.8903830880996
$ws.Range("F20").Value = 32.29005976310985
$ws.Range("G20").Value = 36.21873491593544
$ws.Range("H20").Value = 14.59719531999458
$ws.Range("J20").Value = 9.744413717773465
$ws.Range("N20").Value = 20.96544799484619
$ws.Range("O20").Value = 23.61489016423105

$ws.Range("C21").Value = 5.712241443374218
$ws.Range("D21").Value = 9.667707436493352
$ws.Range("E21").Value = 13.95419934077855
$ws.Range("F21").Value = 32.84002960602491
$ws.Range("G21").Value = 37.33914387723863
$ws.Range("H21").Value = 14.69821609815435
$ws.Range("J21").Value = 9.755725095049558
$ws.Range("N21").Value = 21.70751365554066
$ws.Range("O21").Value = 23.9777698935494

$ws.Range("C22").Value = 5.851889344498632
$ws.Range("D22").Value = 9.671780550156853
$ws.Range("E22").Value = 13.99835171587834
$ws.Range("F22").Value = 33.20197221035591
$ws.Range("G22").Value = 38.0620167835626
$ws.Range("H22").Value = 14.7667327204538
$ws.Range("J22").Value = 9.764635163502508
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 24.21731728487552

$ws.Range("C23").Value = 5.777670146970684
$ws.Range("D23").Value = 9.669485582652163
$ws.Range("E23").Value = 13.97456568595166
$ws.Range("F23").Value = 33.00860570262229
$ws.Range("G23").Value = 37.67714822793447
$ws.Range("H23").Value = 14.72994077576615
$ws.Range("J23").Value = 9.759741260557734
$ws.Range("N23").Value = 21.92877110912573
$ws.Range("O23").Value = 24.08927243296245

$ws.Range("C24").Value = 5.489936827226254
$ws.Range("D24").Value = 9.663458301277922
$ws.Range("E24").Value = 13.88938188403093
$ws.Range("F24").Value = 32.28112168827893
$ws.Range("G24").Value = 36.20028347563829
$ws.Range("H24").Value = 14.59558734499342
$ws.Range("J24").Value = 9.744254543125468
$ws.Range("N24").Value = 20.95310750188672
$ws.Range("O24").Value = 23.60900476706787

$ws.Range("C25").Value = 5.167123255842959
$ws.Range("D25").Value = 9.662676225994339
$ws.Range("E25").Value = 13.80823951132426
$ws.Range("F25").Value = 31.51106847167192
$ws.Range("G25").Value = 34.57607532276501
$ws.Range("H25").Value = 14.46183719371586
$ws.Range("J25").Value = 9.734100479248383
$ws.Range("N25").Value = 19.84905939529495
$ws.Range("O25").Value = 23.10364695539359
